# NIT-9007788941.xlsx — "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Removes the older period-2507 worker rows (keeping one representative row
# per worker) and the duplicate period-2508 block, leaving four worker rows
# that are re-labelled with the new period 2509. Also updates the summary
# totals (Valor Mora, Cant. Trabajadores, Cant. Periodos) to match the
# reduced data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the rows that are no longer needed -----------------------------
# Original data block (rows 16-27) holds 8 workers @ period 2507 (rows16-23)
# and 4 of those same workers repeated @ period 2508 (rows 24-27). We keep
# only one row per worker (Jorge, Rafael, Michel, Jeimy) and drop the rest,
# deleting bottom-to-top so row numbers of not-yet-deleted rows stay stable.
$ws.Rows("27:27").Delete()
$ws.Rows("26:26").Delete()
$ws.Rows("25:25").Delete()
$ws.Rows("24:24").Delete()
$ws.Rows("22:22").Delete()
$ws.Rows("20:20").Delete()
$ws.Rows("18:18").Delete()
$ws.Rows("17:17").Delete()

# --- Relabel the remaining four worker rows with the new period ------------
$ws.Range("E16").Value = "2509"
$ws.Range("E17").Value = "2509"
$ws.Range("E18").Value = "2509"
$ws.Range("E19").Value = "2509"

# Center the "Periodo Mora" column for the data rows (visual tweak that came
# together with the data refresh).
$ws.Range("E16:E19").HorizontalAlignment = -4108

# --- Update the summary figures --------------------------------------------
$ws.Range("E11").Value = 235288
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 1
